$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F1").Value = "Achieved"
$ws.Range("F2").Value = "Dr"
$ws.Range("F4").Value = "Bachelor"

$ws.Range("F2").Select()
